$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1 = "visit", styled like the other header cells (A1:E1)
$ws.Range("F1").Value = "visit"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# New data column F2:F41
$values = @(0,1,0,0,0,0,0,1,0,1,0,0,0,0,0,1,1,0,0,1,0,0,0,0,0,1,0,0,0,1,0,0,0,0,0,0,0,0,1,0)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
}
